# NIT-9011266848.xlsx update
# - Adds a new worker (YURI EIDY VILLADA HOYOS, CC 1129517703) with two
#   overdue periods (2103, 2102) to the "Estado de Cuenta" detail table.
# - Re-orders the existing worker's (ISOLINA MARIA MARVAL PEREZ) period rows
#   so the "Valor Mora" figures line up with period 2407/2406/2405.
# - Refreshes the totals (Valor Mora, Cant. Trabajadores, Cant. Periodos).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Make room for two more detail rows. The table currently ends at row
#    18 (last row carries the bold "closing" border); insert two blank
#    rows right after it so the trailing "firma" block (old rows 23-24)
#    shifts down to rows 25-26, matching the new layout.
# ---------------------------------------------------------------------
$ws.Rows("19:20").Insert()

# Capture the "closing" (bold/bottom-border) row format before it is
# overwritten, and stamp it onto the new last row (20).
$ws.Range("B18:J18").Copy() | Out-Null
$ws.Range("B20:J20").PasteSpecial(-4122) | Out-Null

# Row 19 becomes an ordinary detail row -> copy the regular row format.
$ws.Range("B17:J17").Copy() | Out-Null
$ws.Range("B19:J19").PasteSpecial(-4122) | Out-Null

# Row 18 is no longer the last row of the table -> demote it to the
# regular (non-bold / no bottom border) row format too.
$ws.Range("B17:J17").Copy() | Out-Null
$ws.Range("B18:J18").PasteSpecial(-4122) | Out-Null

$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2. Rewrite the detail rows (16-20) with the refreshed data set.
# ---------------------------------------------------------------------
# Row 16: ISOLINA - period 2407
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "2000021510"
$ws.Range("D16").Value = "ISOLINA MARIA MARVAL PEREZ"
$ws.Range("E16").Value = "2407"
$ws.Range("F16").Value = 60107
$ws.Range("G16").Value = 9016000

# Row 17: ISOLINA - period 2406
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "2000021510"
$ws.Range("D17").Value = "ISOLINA MARIA MARVAL PEREZ"
$ws.Range("E17").Value = "2406"
$ws.Range("F17").Value = 360640
$ws.Range("G17").Value = 9016000

# Row 18: ISOLINA - period 2405
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "2000021510"
$ws.Range("D18").Value = "ISOLINA MARIA MARVAL PEREZ"
$ws.Range("E18").Value = "2405"
$ws.Range("F18").Value = 360640
$ws.Range("G18").Value = 9016000

# Row 19: YURI - period 2103 (new worker)
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1129517703"
$ws.Range("D19").Value = "YURI EIDY VILLADA HOYOS"
$ws.Range("E19").Value = "2103"
$ws.Range("F19").Value = 48413
$ws.Range("G19").Value = 1210308

# Row 20: YURI - period 2102 (new worker, closing row)
$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "1129517703"
$ws.Range("D20").Value = "YURI EIDY VILLADA HOYOS"
$ws.Range("E20").Value = "2102"
$ws.Range("F20").Value = 48413
$ws.Range("G20").Value = 1210308

# ---------------------------------------------------------------------
# 3. Refresh the summary fields above the table.
# ---------------------------------------------------------------------
$ws.Range("E11").Value = 878213
$ws.Range("C13").Value = 2
$ws.Range("F13").Value = 5
